$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the TDL_EC_49_12S row (row 2) - decontamination moved upstream
$ws.Rows.Item(2).Delete()

# Write updated values (recalculated after upstream decontamination)
$ws.Range("A2").Value = "TDL_FS_10_12S"
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 0.65
$ws.Range("D2").Value = 0.314

$ws.Range("A3").Value = "TDL_FS_11_12S"
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 0.843
$ws.Range("D3").Value = 0.482

$ws.Range("A4").Value = "TDL_FS_12_12S"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 0.478
$ws.Range("D4").Value = 0.235

$ws.Range("A5").Value = "TDL_FS_14_12S"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 0.988
$ws.Range("D5").Value = 0.458

$ws.Range("A6").Value = "TDL_FS_15_12S"
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 0.894
$ws.Range("D6").Value = 0.475

$ws.Range("A7").Value = "TDL_FS_16_12S"
$ws.Range("B7").Value = 13
$ws.Range("C7").Value = 0.551
$ws.Range("D7").Value = 0.276

$ws.Range("A8").Value = "TDL_FS_17_12S"
$ws.Range("B8").Value = 13
$ws.Range("C8").Value = 0.147
$ws.Range("D8").Value = 0.044

$ws.Range("A9").Value = "TDL_FS_18_12S"
$ws.Range("B9").Value = 16
$ws.Range("C9").Value = 0.603
$ws.Range("D9").Value = 0.211

$ws.Range("A10").Value = "TDL_FS_19_12S"
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 1.934
$ws.Range("D10").Value = 0.769

$ws.Range("A11").Value = "TDL_FS_1_12S"
$ws.Range("B11").Value = 18
$ws.Range("C11").Value = 1.465
$ws.Range("D11").Value = 0.641

$ws.Range("A12").Value = "TDL_FS_20_12S"
$ws.Range("B12").Value = 20
$ws.Range("C12").Value = 2.162
$ws.Range("D12").Value = 0.811

$ws.Range("A13").Value = "TDL_FS_21_12S"
$ws.Range("B13").Value = 23
$ws.Range("C13").Value = 2.315
$ws.Range("D13").Value = 0.857

$ws.Range("A14").Value = "TDL_FS_22_12S"
$ws.Range("B14").Value = 16
$ws.Range("C14").Value = 1.837
$ws.Range("D14").Value = 0.746

$ws.Range("A15").Value = "TDL_FS_23_12S"
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 2.124
$ws.Range("D15").Value = 0.839

$ws.Range("A16").Value = "TDL_FS_24_12S"
$ws.Range("B16").Value = 22
$ws.Range("C16").Value = 1.863
$ws.Range("D16").Value = 0.786

$ws.Range("A17").Value = "TDL_FS_25_12S"
$ws.Range("B17").Value = 22
$ws.Range("C17").Value = 1.666
$ws.Range("D17").Value = 0.735

$ws.Range("A18").Value = "TDL_FS_26_12S"
$ws.Range("B18").Value = 26
$ws.Range("C18").Value = 1.386
$ws.Range("D18").Value = 0.633

$ws.Range("A19").Value = "TDL_FS_27_12S"
$ws.Range("B19").Value = 22
$ws.Range("C19").Value = 1.544
$ws.Range("D19").Value = 0.709

$ws.Range("A20").Value = "TDL_FS_28_12S"
$ws.Range("B20").Value = 32
$ws.Range("C20").Value = 1.86
$ws.Range("D20").Value = 0.784

$ws.Range("A21").Value = "TDL_FS_29_12S"
$ws.Range("B21").Value = 18
$ws.Range("C21").Value = 1.724
$ws.Range("D21").Value = 0.746

$ws.Range("A22").Value = "TDL_FS_2_12S"
$ws.Range("B22").Value = 17
$ws.Range("C22").Value = 1.553
$ws.Range("D22").Value = 0.709

$ws.Range("A23").Value = "TDL_FS_30_12S"
$ws.Range("B23").Value = 22
$ws.Range("C23").Value = 1.669
$ws.Range("D23").Value = 0.757

$ws.Range("A24").Value = "TDL_FS_31_12S"
$ws.Range("B24").Value = 27
$ws.Range("C24").Value = 2.104
$ws.Range("D24").Value = 0.815

$ws.Range("A25").Value = "TDL_FS_32_12S"
$ws.Range("B25").Value = 22
$ws.Range("C25").Value = 0.773
$ws.Range("D25").Value = 0.294

$ws.Range("A26").Value = "TDL_FS_33_12S"
$ws.Range("B26").Value = 26
$ws.Range("C26").Value = 1.888
$ws.Range("D26").Value = 0.722

$ws.Range("A27").Value = "TDL_FS_34_12S"
$ws.Range("B27").Value = 24
$ws.Range("C27").Value = 1.714
$ws.Range("D27").Value = 0.728

$ws.Range("A28").Value = "TDL_FS_35_12S"
$ws.Range("B28").Value = 29
$ws.Range("C28").Value = 1.25
$ws.Range("D28").Value = 0.544

$ws.Range("A29").Value = "TDL_FS_36_12S"
$ws.Range("B29").Value = 22
$ws.Range("C29").Value = 1.613
$ws.Range("D29").Value = 0.74

$ws.Range("A30").Value = "TDL_FS_37_12S"
$ws.Range("B30").Value = 26
$ws.Range("C30").Value = 1.938
$ws.Range("D30").Value = 0.738

$ws.Range("A31").Value = "TDL_FS_38_12S"
$ws.Range("B31").Value = 29
$ws.Range("C31").Value = 2.034
$ws.Range("D31").Value = 0.751

$ws.Range("A32").Value = "TDL_FS_39_12S"
$ws.Range("B32").Value = 43
$ws.Range("C32").Value = 1.989
$ws.Range("D32").Value = 0.744

$ws.Range("A33").Value = "TDL_FS_3_12S"
$ws.Range("B33").Value = 16
$ws.Range("C33").Value = 1.182
$ws.Range("D33").Value = 0.544

$ws.Range("A34").Value = "TDL_FS_40_12S"
$ws.Range("B34").Value = 23
$ws.Range("C34").Value = 1.938
$ws.Range("D34").Value = 0.785

$ws.Range("A35").Value = "TDL_FS_41_12S"
$ws.Range("B35").Value = 21
$ws.Range("C35").Value = 1.68
$ws.Range("D35").Value = 0.719

$ws.Range("A36").Value = "TDL_FS_42_12S"
$ws.Range("B36").Value = 26
$ws.Range("C36").Value = 2.077
$ws.Range("D36").Value = 0.8

$ws.Range("A37").Value = "TDL_FS_43_12S"
$ws.Range("B37").Value = 20
$ws.Range("C37").Value = 1.853
$ws.Range("D37").Value = 0.759

$ws.Range("A38").Value = "TDL_FS_44_12S"
$ws.Range("B38").Value = 21
$ws.Range("C38").Value = 1.706
$ws.Range("D38").Value = 0.723

$ws.Range("A39").Value = "TDL_FS_45_12S"
$ws.Range("B39").Value = 25
$ws.Range("C39").Value = 1.715
$ws.Range("D39").Value = 0.732

$ws.Range("A40").Value = "TDL_FS_46_12S"
$ws.Range("B40").Value = 19
$ws.Range("C40").Value = 1.299
$ws.Range("D40").Value = 0.62

$ws.Range("A41").Value = "TDL_FS_47_12S"
$ws.Range("B41").Value = 43
$ws.Range("C41").Value = 2.306
$ws.Range("D41").Value = 0.779

$ws.Range("A42").Value = "TDL_FS_48_12S"
$ws.Range("B42").Value = 39
$ws.Range("C42").Value = 2.612
$ws.Range("D42").Value = 0.885

$ws.Range("A43").Value = "TDL_FS_4_12S"
$ws.Range("B43").Value = 13
$ws.Range("C43").Value = 1.118
$ws.Range("D43").Value = 0.511

$ws.Range("A44").Value = "TDL_FS_5_12S"
$ws.Range("B44").Value = 16
$ws.Range("C44").Value = 1.014
$ws.Range("D44").Value = 0.435

$ws.Range("A45").Value = "TDL_FS_6_12S"
$ws.Range("B45").Value = 12
$ws.Range("C45").Value = 0.432
$ws.Range("D45").Value = 0.178

$ws.Range("A46").Value = "TDL_FS_7_12S"
$ws.Range("B46").Value = 15
$ws.Range("C46").Value = 0.534
$ws.Range("D46").Value = 0.24

$ws.Range("A47").Value = "TDL_FS_8_12S"
$ws.Range("B47").Value = 19
$ws.Range("C47").Value = 0.851
$ws.Range("D47").Value = 0.527

$ws.Range("A48").Value = "TDL_FS_9_12S"
$ws.Range("B48").Value = 16
$ws.Range("C48").Value = 1.256
$ws.Range("D48").Value = 0.671
